$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 13.11017566666667
$ws.Range("H2").Value = 39.330527
$ws.Range("I2").Value = 0.1657114824704502
$ws.Range("J2").Value = 0.1657114824704501
$ws.Range("M2").Value = 0.9386610000000001
$ws.Range("N2").Value = 2.815983
$ws.Range("O2").Value = 0.04181245246793033
$ws.Range("P2").Value = 0.04181245246793032
$ws.Range("Q2").Value = 12.306010601449
$ws.Range("R2").Value = 110.754095413041
$ws.Range("S2").Value = 0.006928803484185967
$ws.Range("T2").Value = 0.006928803484185965

# Row 3
$ws.Range("G3").Value = 13.11017566666667
$ws.Range("H3").Value = 39.330527
$ws.Range("I3").Value = 0.1657114824704502
$ws.Range("J3").Value = 0.1657114824704501
$ws.Range("O3").Value = 0.1106393125456779
$ws.Range("P3").Value = 0.1106393125456779
$ws.Range("Q3").Value = 32.56275278682634
$ws.Range("R3").Value = 293.064775081437
$ws.Range("S3").Value = 0.01833420450145576
$ws.Range("T3").Value = 0.01833420450145576

# Row 4
$ws.Range("G4").Value = 13.11017566666667
$ws.Range("H4").Value = 39.330527
$ws.Range("I4").Value = 0.1657114824704502
$ws.Range("J4").Value = 0.1657114824704501
$ws.Range("O4").Value = 0.8475482349863918
$ws.Range("P4").Value = 0.8475482349863918
$ws.Range("Q4").Value = 249.445726078411
$ws.Range("R4").Value = 2245.011534705699
$ws.Range("S4").Value = 0.1404484744848084
$ws.Range("T4").Value = 0.1404484744848084

# Row 5
$ws.Range("G5").Value = 51.42568199999999
$ws.Range("I5").Value = 0.6500161567583834
$ws.Range("J5").Value = 0.6500161567583833
$ws.Range("M5").Value = 0.9386610000000001
$ws.Range("N5").Value = 2.815983
$ws.Range("O5").Value = 0.04181245246793033
$ws.Range("P5").Value = 0.04181245246793032
$ws.Range("Q5").Value = 48.271282091802
$ws.Range("R5").Value = 434.441538826218
$ws.Range("S5").Value = 0.02717876965784665
$ws.Range("T5").Value = 0.02717876965784665

# Row 6
$ws.Range("G6").Value = 51.42568199999999
$ws.Range("I6").Value = 0.6500161567583834
$ws.Range("J6").Value = 0.6500161567583833
$ws.Range("O6").Value = 0.1106393125456779
$ws.Range("P6").Value = 0.1106393125456779
$ws.Range("S6").Value = 0.07191734072733115
$ws.Range("T6").Value = 0.07191734072733112

# Row 7
$ws.Range("G7").Value = 51.42568199999999
$ws.Range("I7").Value = 0.6500161567583834
$ws.Range("J7").Value = 0.6500161567583833
$ws.Range("O7").Value = 0.8475482349863918
$ws.Range("P7").Value = 0.8475482349863918
$ws.Range("Q7").Value = 978.4702289064778
$ws.Range("R7").Value = 8806.2320601583
$ws.Range("S7").Value = 0.5509200463732056
$ws.Range("T7").Value = 0.5509200463732056

# Row 8
$ws.Range("I8").Value = 0.1842723607711665
$ws.Range("J8").Value = 0.1842723607711665
$ws.Range("M8").Value = 0.9386610000000001
$ws.Range("N8").Value = 2.815983
$ws.Range("O8").Value = 0.04181245246793033
$ws.Range("P8").Value = 0.04181245246793032
$ws.Range("Q8").Value = 13.684372328323
$ws.Range("R8").Value = 123.159350954907
$ws.Range("S8").Value = 0.007704879325897709
$ws.Range("T8").Value = 0.007704879325897707

# Row 9
$ws.Range("I9").Value = 0.1842723607711665
$ws.Range("J9").Value = 0.1842723607711665
$ws.Range("O9").Value = 0.1106393125456779
$ws.Range("P9").Value = 0.1106393125456779
$ws.Range("S9").Value = 0.02038776731689101
$ws.Range("T9").Value = 0.020387767316891

# Row 10
$ws.Range("I10").Value = 0.1842723607711665
$ws.Range("J10").Value = 0.1842723607711665
$ws.Range("O10").Value = 0.8475482349863918
$ws.Range("P10").Value = 0.8475482349863918
$ws.Range("S10").Value = 0.1561797141283778
$ws.Range("T10").Value = 0.1561797141283778
